# Update gh-pages to output generated at 456a3b4
# This applies updated "want to go" counts (column F) across the four sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F7").Value = 2568
$wsExhibition.Range("F8").Value = 949
$wsExhibition.Range("F9").Value = 18943
$wsExhibition.Range("F10").Value = 60
$wsExhibition.Range("F11").Value = 2017
$wsExhibition.Range("F14").Value = 356
$wsExhibition.Range("F15").Value = 622
$wsExhibition.Range("F16").Value = 203
$wsExhibition.Range("F20").Value = 50
$wsExhibition.Range("F21").Value = 219

# Sheet "演出" (Performances)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F10").Value = 239
$wsPerformance.Range("F11").Value = 239
$wsPerformance.Range("F14").Value = 10
$wsPerformance.Range("F19").Value = 28

# Sheet "本地生活" (Local Life)
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F2").Value = 5920

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5920
$wsAll.Range("F14").Value = 2568
$wsAll.Range("F15").Value = 949
$wsAll.Range("F16").Value = 18944
$wsAll.Range("F19").Value = 60
$wsAll.Range("F21").Value = 239
$wsAll.Range("F22").Value = 239
$wsAll.Range("F23").Value = 2017
$wsAll.Range("F26").Value = 356
$wsAll.Range("F27").Value = 622
$wsAll.Range("F28").Value = 203
$wsAll.Range("F31").Value = 10
$wsAll.Range("F34").Value = 50
$wsAll.Range("F36").Value = 219
$wsAll.Range("F40").Value = 28
